$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11 ("Export to PDF") - content placeholder gets a position/size and
# three new bullet paragraphs.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$body11 = $s11.Shapes.Item(2)

$body11.Left = 89.875
$body11.Top = 177.125
$body11.Width = 779.9999389648438
$body11.Height = 234.0055389404297

$tr11 = $body11.TextFrame.TextRange
$tr11.Text = "Using open source iText library for java"
[void]$tr11.InsertAfter([char]13 + "Using rectangle to implement block control for positioning paragraph text")
[void]$tr11.InsertAfter([char]13 + "Using table to display order items of invoices")

# Force a clean run split around "iText" (no formatting change, just a
# dedicated run boundary matching the authored markup).
$body11.TextFrame.TextRange.Characters(19, 5).Text = "iText"

# ---------------------------------------------------------------------------
# Slide 9 ("DAO class" -> "DAO (DATABASE class)") - title text/position and
# content placeholder position/text.
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1)

$title9.Left = 58.57071304321289
$title9.Top = 102.2608871459961

$title9.TextFrame.TextRange.Text = "DAO (DATABASE class)"

$body9 = $s9.Shapes.Item(2)

$body9.Left = 58.5706672668457
$body9.Top = 201.391357421875
$body9.Width = 779.9999389648438
$body9.Height = 302.6086730957031

$tr9 = $body9.TextFrame.TextRange
$tr9.Text = 'Using database connection transaction ensure the data integrity for issuing invoice (insert record into invoices table, meanwhile modify status to  "complete" for orders through update orders table)'
[void]$tr9.InsertAfter([char]13 + "Overload getOrders method for adapting different search functionality")
[void]$tr9.InsertAfter([char]13 + "Using StringBulider instead of String for concatenation sql statements, which are easier readable and more elegant")

# Force clean run splits matching the authored markup (spell-check style
# isolation of "getOrders", "StringBulider" and "sql").
$body9.TextFrame.TextRange.Characters(209, 9).Text = "getOrders"
$body9.TextFrame.TextRange.Characters(276, 13).Text = "StringBulider"
$body9.TextFrame.TextRange.Characters(326, 3).Text = "sql"
